$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Overview sheet: status cells move from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both languages/rows, and the
#    zh-cn / de-de status columns (E, F) are widened to fit the new text.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus

# ColumnWidth (character units) = stored sheet width - 5/6
$wsOverview.Range("E:F").ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------------
# Helper values shared by the two language sheets.
# ---------------------------------------------------------------------------
$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/149bcc218bd00c6efe143eda8dcc0568e4b1bcb9/e2e/37293e44-a6c0-4225-9a0d-4fcba3b5eb01.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/149bcc218bd00c6efe143eda8dcc0568e4b1bcb9/e2e/d86bf385-42c1-476a-8eb2-d5c78d48af64.md"
$dispA = "37293e44-a6c0-4225-9a0d-4fcba3b5eb01.md"
$dispB = "d86bf385-42c1-476a-8eb2-d5c78d48af64.md"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: the handback run now has a target file (column I) and a
#    handback file (column J) for each of the two rows. Rebuild the
#    hyperlinks collection so the new ones line up in document order
#    (A2, I2, A3, I3) and keep the matching "HyperLink" style.
# ---------------------------------------------------------------------------
$wsZhCn.Range("C:C").ColumnWidth = 29.144371396019366
$wsZhCn.Range("I:I").ColumnWidth = 39.166666666666664
$wsZhCn.Range("J:J").ColumnWidth = 39.166666666666664

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlA, "", "", $dispA)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlA, "", "", $dispA)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlB, "", "", $dispB)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlB, "", "", $dispB)

$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276

$wsZhCn.Range("I2").Value2 = $dispA
$wsZhCn.Range("J2").Value2 = "37293e44-a6c0-4225-9a0d-4fcba3b5eb01.66a7b4bc72de4845fbd5886e92f5eb799ef1e88b.zh-cn.xlf"
$wsZhCn.Range("I3").Value2 = $dispB
$wsZhCn.Range("J3").Value2 = "d86bf385-42c1-476a-8eb2-d5c78d48af64.717192b6f2e9029d8fff8e5a23389ae928ca68fc.zh-cn.xlf"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of change as zh-cn, plus the handback datetime
#    (column K) moves from the "never handed back" placeholder to the real
#    handback timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Range("C:C").ColumnWidth = 29.144371396019366
$wsDeDe.Range("I:I").ColumnWidth = 39.166666666666664
$wsDeDe.Range("J:J").ColumnWidth = 39.166666666666664

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlA, "", "", $dispA)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlA, "", "", $dispA)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlB, "", "", $dispB)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlB, "", "", $dispB)

$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276

$wsDeDe.Range("I2").Value2 = $dispA
$wsDeDe.Range("J2").Value2 = "37293e44-a6c0-4225-9a0d-4fcba3b5eb01.66a7b4bc72de4845fbd5886e92f5eb799ef1e88b.de-de.xlf"
$wsDeDe.Range("K2").Value2 = "2016-09-03 16:54:46"
$wsDeDe.Range("I3").Value2 = $dispB
$wsDeDe.Range("J3").Value2 = "d86bf385-42c1-476a-8eb2-d5c78d48af64.717192b6f2e9029d8fff8e5a23389ae928ca68fc.de-de.xlf"
$wsDeDe.Range("K3").Value2 = "2016-09-03 16:54:46"

Write-Host "Handback report generated."
